$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply "custom accuracy": round the last data row's measurements (row 5,
# columns B:AH) down to 2 decimal places instead of 3.
$ws.Range("B5").Value = 8.17
$ws.Range("C5").Value = 5.8
$ws.Range("D5").Value = 0.83
$ws.Range("E5").Value = 17.46
$ws.Range("F5").Value = 14.45
$ws.Range("G5").Value = 6.43
$ws.Range("H5").Value = 26.75
$ws.Range("I5").Value = 9.89
$ws.Range("J5").Value = 4.29
$ws.Range("K5").Value = 6.44
$ws.Range("L5").Value = 7.09
$ws.Range("M5").Value = 7.32
$ws.Range("N5").Value = 2.05
$ws.Range("O5").Value = 6.39
$ws.Range("P5").Value = 9.02
$ws.Range("Q5").Value = 5.52
$ws.Range("R5").Value = 0.76
$ws.Range("S5").Value = 0.46
$ws.Range("T5").Value = 89.96
$ws.Range("U5").Value = 17.94
$ws.Range("V5").Value = 5.9
$ws.Range("W5").Value = 11.9
$ws.Range("X5").Value = 6.42
$ws.Range("Y5").Value = 0.74
$ws.Range("Z5").Value = 12.6
$ws.Range("AA5").Value = 5.21
$ws.Range("AB5").Value = 4.72
$ws.Range("AC5").Value = 5.53
$ws.Range("AD5").Value = 7.38
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 24.18
$ws.Range("AG5").Value = 3.26
$ws.Range("AH5").Value = 7.38

# Drop the last data row (row 6) entirely -- the "1000 data points" trim
# mentioned in the commit message -- which also shrinks the used range from
# A1:AH6 down to A1:AH5.
$ws.Rows.Item(6).Delete()
